$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header
$ws.Range("G1").Value = "Khóa học"

# Fill in the values for the new column
$ws.Range("G2").Value = "2019-2023"
$ws.Range("G3").Value = "2019-2023"
$ws.Range("G4").Value = "2019-2023"
$ws.Range("G5").Value = "2019-2023"

# Set column width for new column G (closest achievable value to 12.6640625)
$ws.Columns.Item(7).ColumnWidth = 11.830729166666666

# Update selection to match new cursor position
$ws.Range("F9").Select()
